$d = $word.ActiveDocument

# The document currently ends with a trailing empty paragraph. We need to
# append, after it:
#   1. a new paragraph containing only a manual page break
#   2. a new paragraph that starts the next page ("This is a new line"),
#      carrying the lastRenderedPageBreak marker Word stamps on the first
#      run after a page break
#   3. a new (empty) paragraph that just holds the "_GoBack" bookmark that
#      Word drops at the last edited location when the file is saved

# --- 1 & 2: create a fresh paragraph mark, then stream in the page break
#     and the "This is a new line" paragraph as real WordprocessingML via
#     InsertXML so the lastRenderedPageBreak marker round-trips exactly.
$endPos = $d.Content.End - 1
$rng = $d.Range($endPos, $endPos)
$rng.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRng = $newPara.Range
$newRng.Collapse(1) | Out-Null

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p><w:r><w:br w:type="page"/></w:r></w:p>' +
       '<w:p><w:r><w:lastRenderedPageBreak/><w:t>This is a new line</w:t></w:r></w:p>' +
       '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRng.InsertXML($xml) | Out-Null

# --- 3: bookmark the trailing empty paragraph that remains after the text
#     we just streamed in, exactly like Word leaves "_GoBack" where the
#     cursor last was on save.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $lastPara.Range) | Out-Null
